$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the text values in E2/E3 (phone numbers) with actual numeric values
$ws.Range("E2").Value = 1234567890
$ws.Range("E3").Value = 2112345678

# Update the active selection on the sheet from F6 to D6
$ws.Range("D6").Select()
